$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "xcv"
$ws.Range("B11").Value = 15
$ws.Range("C11").Value = "g"
$ws.Range("D11").Value = 15

$ws.Range("D11").Select()
